$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "-"
$ws.Range("B4").Value = "-"
$ws.Range("E4").Value = "-"
$ws.Range("F4").Value = "-"
$ws.Range("E6").Value = "-"
$ws.Range("F6").Value = "-"

$ws.Range("B10").Value = "[Aline S. M.-Metalografia-2A, -, -, Aline S. M.-Metalografia-2A]"

$ws.Range("B11").Value = "[Aline S. M.-Metalografia-2A, Andre B.-Elet. Digi. Básica-2A, Andre B.-Elet. Digi. Básica-2A, Aline S. M.-Metalografia-2A]"
$ws.Range("C11").Value = "Ludoff-Maq. Term. Fluxo"

$ws.Range("B12").Value = "[Aderci-Tornearia-2A, Andre B.-Elet. Digi. Básica-2A, Andre B.-Elet. Digi. Básica-2A, Aderci-Tornearia-2A]"
$ws.Range("C12").Value = "Ludoff-Maq. Term. Fluxo"
$ws.Range("D12").Value = "Gilberto-Mec. Tec. Res. Mat."

$ws.Range("B14").Value = "[Aderci-Tornearia-2A, Claudinei-Des. Maq. CAD-T2-2A, Claudinei-Des. Maq. CAD-T2-2A, Aderci-Tornearia-2A]"
$ws.Range("D14").Value = "Gilberto-Mec. Tec. Res. Mat."

$ws.Range("B15").Value = "[Suzanny-Des. Maq. CAD-T1-2A, Claudinei-Des. Maq. CAD-T2-2A, Suzanny-Des. Maq. CAD-T1-2A, Suzanny-Des. Maq. CAD-T1-2A]"
$ws.Range("F15").Value = "[Anderson-Ajustagem-2A, Anderson-Ajustagem-2A, Anderson-Ajustagem-2A, Anderson-Ajustagem-2A]"

$ws.Range("B18").Value = "-"
$ws.Range("C18").Value = "-"
$ws.Range("D18").Value = "-"
$ws.Range("F18").Value = "-"

$ws.Range("B19").Value = "-"
$ws.Range("C19").Value = "-"
$ws.Range("D19").Value = "-"
$ws.Range("F19").Value = "-"

$ws.Range("B20").Value = "-"
$ws.Range("C20").Value = "-"
$ws.Range("D20").Value = "-"
$ws.Range("F20").Value = "-"

$ws.Range("B21").Value = "-"
$ws.Range("C21").Value = "-"
$ws.Range("D21").Value = "-"
$ws.Range("F21").Value = "-"
